# Weekly update: insert this week's price record for
# "Agrícola del Norte S.A. de Arica - Acelga" as a new row right above
# last week's row (old row 68), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 68 - everything currently at/after
# row 68 (including styles) shifts down to row+1.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with this week's data.
$ws.Range("A68").Value = 1
$ws.Range("B68").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C68").Value = "Arica y Parinacota"
$ws.Range("D68").Value = 44889
$ws.Range("E68").Value = 15
$ws.Range("F68").Value = 100112009
$ws.Range("G68").Value = "Acelga"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 500
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = 1750
$ws.Range("N68").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 583
$ws.Range("Q68").Value = 3
$ws.Range("R68").Value = "Hortaliza"
